$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 32 values
$ws.Range("D32").Value = 65
$ws.Range("E32").Value = 70
$ws.Range("F32").Value = 75

# Update existing row 33 values
$ws.Range("C33").Value = 60
$ws.Range("D33").Value = 70
$ws.Range("F33").Value = 90

# Update existing row 34 values
$ws.Range("C34").Value = 85

# Add new rows 40-42
$ws.Range("A40").Value = "low_speed"
$ws.Range("B40").Value = "trapezoidal_mf"
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 50
$ws.Range("E40").Value = 50
$ws.Range("F40").Value = 100

$ws.Range("A41").Value = "mid_speed"
$ws.Range("B41").Value = "trapezoidal_mf"
$ws.Range("C41").Value = 50
$ws.Range("D41").Value = 100
$ws.Range("E41").Value = 100
$ws.Range("F41").Value = 150

$ws.Range("A42").Value = "high_speed"
$ws.Range("B42").Value = "trapezoidal_mf"
$ws.Range("C42").Value = 100
$ws.Range("D42").Value = 150
$ws.Range("E42").Value = 150
$ws.Range("F42").Value = 200

# Update the sheet view (scroll position / selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E40").Select()
